$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.158107333333333
$ws.Range("H2").Value = 6.474322
$ws.Range("I2").Value = 0.04789414999021155
$ws.Range("J2").Value = 0.04789414999021156
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.806900666666666
$ws.Range("N2").Value = 26.420702
$ws.Range("O2").Value = 0.1733678197953833
$ws.Range("P2").Value = 0.1733678197953834
$ws.Range("Q2").Value = 19.00623691267155
$ws.Range("R2").Value = 171.056132214044
$ws.Range("S2").Value = 0.008303304364756056
$ws.Range("T2").Value = 0.00830330436475606

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.158107333333333
$ws.Range("H3").Value = 6.474322
$ws.Range("I3").Value = 0.04789414999021155
$ws.Range("J3").Value = 0.04789414999021156
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.76689066666667
$ws.Range("N3").Value = 56.30067200000001
$ws.Range("O3").Value = 0.3694347242421866
$ws.Range("P3").Value = 0.3694347242421866
$ws.Range("Q3").Value = 40.50096437159822
$ws.Range("R3").Value = 364.508679344384
$ws.Range("S3").Value = 0.01769376209444773
$ws.Range("T3").Value = 0.01769376209444773

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.158107333333333
$ws.Range("H4").Value = 6.474322
$ws.Range("I4").Value = 0.04789414999021155
$ws.Range("J4").Value = 0.04789414999021156
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.225144
$ws.Range("N4").Value = 69.675432
$ws.Range("O4").Value = 0.4571974559624301
$ws.Range("P4").Value = 0.4571974559624301
$ws.Range("Q4").Value = 50.12235358412266
$ws.Range("R4").Value = 451.101182257104
$ws.Range("S4").Value = 0.02189708353100777
$ws.Range("T4").Value = 0.02189708353100777

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 41.35786966666667
$ws.Range("H5").Value = 124.073609
$ws.Range("I5").Value = 0.9178412873614971
$ws.Range("J5").Value = 0.9178412873614971
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.806900666666666
$ws.Range("N5").Value = 26.420702
$ws.Range("O5").Value = 0.1733678197953833
$ws.Range("P5").Value = 0.1733678197953834
$ws.Range("Q5").Value = 364.2346499392797
$ws.Range("R5").Value = 3278.111849453518
$ws.Range("S5").Value = 0.1591241429080507
$ws.Range("T5").Value = 0.1591241429080507

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 41.35786966666667
$ws.Range("H6").Value = 124.073609
$ws.Range("I6").Value = 0.9178412873614971
$ws.Range("J6").Value = 0.9178412873614971
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.76689066666667
$ws.Range("N6").Value = 56.30067200000001
$ws.Range("O6").Value = 0.3694347242421866
$ws.Range("P6").Value = 0.3694347242421866
$ws.Range("Q6").Value = 776.1586182405832
$ws.Range("R6").Value = 6985.427564165249
$ws.Range("S6").Value = 0.3390824428944882
$ws.Range("T6").Value = 0.3390824428944882

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 41.35786966666667
$ws.Range("H7").Value = 124.073609
$ws.Range("I7").Value = 0.9178412873614971
$ws.Range("J7").Value = 0.9178412873614971
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.225144
$ws.Range("N7").Value = 69.675432
$ws.Range("O7").Value = 0.4571974559624301
$ws.Range("P7").Value = 0.4571974559624301
$ws.Range("Q7").Value = 960.5424785415653
$ws.Range("R7").Value = 8644.882306874088
$ws.Range("S7").Value = 0.4196347015589582
$ws.Range("T7").Value = 0.4196347015589582

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.543959
$ws.Range("H8").Value = 4.631876999999999
$ws.Range("I8").Value = 0.03426456264829138
$ws.Range("J8").Value = 0.03426456264829137
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.806900666666666
$ws.Range("N8").Value = 26.420702
$ws.Range("O8").Value = 0.1733678197953833
$ws.Range("P8").Value = 0.1733678197953834
$ws.Range("Q8").Value = 13.597493546406
$ws.Range("R8").Value = 122.377441917654
$ws.Range("S8").Value = 0.005940372522576603
$ws.Range("T8").Value = 0.005940372522576603

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.543959
$ws.Range("H9").Value = 4.631876999999999
$ws.Range("I9").Value = 0.03426456264829138
$ws.Range("J9").Value = 0.03426456264829137
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.76689066666667
$ws.Range("N9").Value = 56.30067200000001
$ws.Range("O9").Value = 0.3694347242421866
$ws.Range("P9").Value = 0.3694347242421866
$ws.Range("Q9").Value = 28.975309746816
$ws.Range("R9").Value = 260.777787721344
$ws.Range("S9").Value = 0.01265851925325065
$ws.Range("T9").Value = 0.01265851925325065

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.543959
$ws.Range("H10").Value = 4.631876999999999
$ws.Range("I10").Value = 0.03426456264829138
$ws.Range("J10").Value = 0.03426456264829137
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.225144
$ws.Range("N10").Value = 69.675432
$ws.Range("O10").Value = 0.4571974559624301
$ws.Range("P10").Value = 0.4571974559624301
$ws.Range("Q10").Value = 35.858670105096
$ws.Range("R10").Value = 322.728030945864
$ws.Range("S10").Value = 0.01566567087246412
$ws.Range("T10").Value = 0.01566567087246412

